# Applies the "Add files via upload" edit:
#  - Slide 2 (subtitle placeholder): shorten the first bullet's text and
#    normalize the second bullet's paragraph spacing / bullet size.
#  - Slide 4 (body placeholder): rename "Create a new notebook" to
#    "Upload a notebook", collapse its instruction bullet to the upload
#    variant, and demote that + the GPU-runtime instruction bullet one
#    indent level.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - "Use your UC email..." subtitle box
# ---------------------------------------------------------------------
$s2   = $p.Slides.Item(2)
$tr2  = $s2.Shapes.Item(1).TextFrame.TextRange

# Paragraph 1: trim the sentence, dropping the "...and then click Slides" tail.
$para1 = $tr2.Paragraphs(1, 1)
$run1  = $para1.Runs(1, 1)
$run1.Text = "Use your UC email address to login your Google account."

# Paragraph 2: "Click https://... and download the entire folder"
$para2 = $tr2.Paragraphs(2, 1)
$pf2   = $para2.ParagraphFormat
$pf2.SpaceBefore = 10
$pf2.SpaceAfter  = 0
$pf2.Bullet.Font.Size = 24
$para2.IndentLevel = 1

# ---------------------------------------------------------------------
# Slide 4 - Colab setup instructions box
# ---------------------------------------------------------------------
$s4  = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(1).TextFrame.TextRange

# Paragraph 3: "Create a new notebook" -> "Upload a notebook"
$para3 = $tr4.Paragraphs(3, 1)
$run3  = $para3.Runs(1, 1)
$run3.Text = "Upload a notebook"

# Paragraph 4: "File -> New Notebook. Either Python 2 or 3 ..." ->
# single run "File -> Upload Notebook. ", demoted one level.
$para4 = $tr4.Paragraphs(4, 1)
$run4b = $para4.Runs(2, 1)
$run4b.Text = ""
$run4a = $para4.Runs(1, 1)
$run4a.Text = "File -> Upload Notebook. "
$para4.IndentLevel = 2

# Paragraph 6: "Runtime -> Change runtime type -> ... Free GPU cycles!"
# stays the same text, just demoted one level.
$para6 = $tr4.Paragraphs(6, 1)
$para6.IndentLevel = 2
